# Fruta / hortaliza, semanal
# Update weekly price data for "Femacal de La Calera - Frambuesa" sheet.
# Rows 2-17 get their Fecha / Volumen / Precio minimo / Precio maximo /
# Precio promedio ponderado / Origen / Precio $/Kg values refreshed, and a
# new data row (18) is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44588

# Row 3
$ws.Range("D3").Value = 44214
$ws.Range("M3").Value = 48
$ws.Range("N3").Value = 6000
$ws.Range("O3").Value = 6000
$ws.Range("P3").Value = 6000
$ws.Range("R3").Value = "Provincia de Linares"
$ws.Range("S3").Value = 3000

# Row 4
$ws.Range("D4").Value = 44586
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 7000
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 7000
$ws.Range("R4").Value = "Provincia de Curicó"
$ws.Range("S4").Value = 3500

# Row 5
$ws.Range("D5").Value = 44592
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = 8000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 8000
$ws.Range("R5").Value = "Provincia de Linares"
$ws.Range("S5").Value = 4000

# Row 6
$ws.Range("D6").Value = 44974
$ws.Range("M6").Value = 130
$ws.Range("N6").Value = 7000
$ws.Range("O6").Value = 7500
$ws.Range("P6").Value = 7269
$ws.Range("S6").Value = 3634

# Row 7
$ws.Range("D7").Value = 44959
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = 7000
$ws.Range("P7").Value = 7000
$ws.Range("S7").Value = 3500

# Row 8
$ws.Range("D8").Value = 44582
$ws.Range("M8").Value = 150
$ws.Range("O8").Value = 6500
$ws.Range("P8").Value = 6233
$ws.Range("R8").Value = "Provincia de Curicó"
$ws.Range("S8").Value = 3116

# Row 9
$ws.Range("D9").Value = 44585
$ws.Range("M9").Value = 160
$ws.Range("N9").Value = 6500
$ws.Range("P9").Value = 6750
$ws.Range("R9").Value = "Provincia de Curicó"
$ws.Range("S9").Value = 3375

# Row 10
$ws.Range("D10").Value = 44614
$ws.Range("M10").Value = 45
$ws.Range("N10").Value = 6000
$ws.Range("O10").Value = 6000
$ws.Range("P10").Value = 6000
$ws.Range("R10").Value = "Provincia de Linares"
$ws.Range("S10").Value = 3000

# Row 11
$ws.Range("D11").Value = 44209
$ws.Range("M11").Value = 58
$ws.Range("N11").Value = 6000
$ws.Range("O11").Value = 6000
$ws.Range("P11").Value = 6000
$ws.Range("S11").Value = 3000

# Row 12
$ws.Range("D12").Value = 44211
$ws.Range("M12").Value = 45
$ws.Range("O12").Value = 6000
$ws.Range("P12").Value = 6000
$ws.Range("S12").Value = 3000

# Row 14
$ws.Range("D14").Value = 44628
$ws.Range("M14").Value = 40
$ws.Range("R14").Value = "Provincia de Linares"

# Row 15
$ws.Range("D15").Value = 44627
$ws.Range("M15").Value = 45

# Row 16
$ws.Range("D16").Value = 44606
$ws.Range("M16").Value = 45
$ws.Range("N16").Value = 7000
$ws.Range("O16").Value = 7000
$ws.Range("P16").Value = 7000
$ws.Range("S16").Value = 3500

# Row 17
$ws.Range("D17").Value = 44960
$ws.Range("M17").Value = 40
$ws.Range("N17").Value = 7000
$ws.Range("O17").Value = 7000
$ws.Range("P17").Value = 7000
$ws.Range("R17").Value = "Provincia de Curicó"
$ws.Range("S17").Value = 3500

# New row 18: append a new weekly record, copying formatting from row 17.
$ws.Range("A17:T17").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").PasteSpecial(-4123)

$ws.Range("A18").Value = 3
$ws.Range("B18").Value = "Femacal de La Calera"
$ws.Range("C18").Value = "Coquimbo"
$ws.Range("D18").Value = 44589
$ws.Range("E18").Value = 5
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100101
$ws.Range("H18").Value = "Berries"
$ws.Range("I18").Value = 100101004
$ws.Range("J18").Value = "Frambuesa"
$ws.Range("K18").Value = "Sin especificar"
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 60
$ws.Range("N18").Value = 6000
$ws.Range("O18").Value = 6000
$ws.Range("P18").Value = 6000
$ws.Range("Q18").Value = "$/bandeja 2 kilos"
$ws.Range("R18").Value = "Provincia de Curicó"
$ws.Range("S18").Value = 3000
$ws.Range("T18").Value = 2
